# Generate Report for Handback
# For each language sheet (zh-cn, de-de), the file 756e0514-...md has come
# back from handback and is now in sync with en-US. Update row 2 (the row
# for that source file) to:
#   - change the Status column (C) from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - fill in "Latest Target File" (F) and "Latest Handback File" (G) with
#     links to the handed-off / handed-back files
#   - set "Latest Handback DateTime" (H) to the real handback timestamp
#     (replacing the "0001-01-01 00:00:00" placeholder)
# Row 3 (cdf7d01a-...md) has not been handed back yet, so it is left as-is.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/10011ddec33d99c4a02e5fef07fbdbe0eb2a7b6d/e2e/756e0514-3dd6-4c71-8197-12542bec0e46.md", [type]::Missing, [type]::Missing, "756e0514-3dd6-4c71-8197-12542bec0e46.md")
$ws.Range("F2").Font.Underline = 2
$ws.Range("F2").Font.Color = 15570276

$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d290ef630d7cb5c046675d051d6fb292da3832bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/756e0514-3dd6-4c71-8197-12542bec0e46.aa50abacea9f47958792410feec0f185133fb1df.zh-cn.xlf", [type]::Missing, [type]::Missing, "756e0514-3dd6-4c71-8197-12542bec0e46.aa50abacea9f47958792410feec0f185133fb1df.zh-cn.xlf")
$ws.Range("G2").Font.Underline = 2
$ws.Range("G2").Font.Color = 15570276

$ws.Range("H2").Value = "2016-03-12 08:26:51"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/10011ddec33d99c4a02e5fef07fbdbe0eb2a7b6d/e2e/756e0514-3dd6-4c71-8197-12542bec0e46.md", [type]::Missing, [type]::Missing, "756e0514-3dd6-4c71-8197-12542bec0e46.md")
$ws.Range("F2").Font.Underline = 2
$ws.Range("F2").Font.Color = 15570276

$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e8db92b401f1e74380e47a3fa1b1b191cf2ae30f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/756e0514-3dd6-4c71-8197-12542bec0e46.aa50abacea9f47958792410feec0f185133fb1df.de-de.xlf", [type]::Missing, [type]::Missing, "756e0514-3dd6-4c71-8197-12542bec0e46.aa50abacea9f47958792410feec0f185133fb1df.de-de.xlf")
$ws.Range("G2").Font.Underline = 2
$ws.Range("G2").Font.Color = 15570276

$ws.Range("H2").Value = "2016-03-12 08:26:57"
